$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with two new columns P1=14, Q1=15,
# matching the style of the existing header cells (bold/centered/bordered).
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)
$ws.Range("Q1").PasteSpecial(-4122)
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15

# Row pattern (columns B..Q, i.e. col index 2..17) that every data row
# (2 through 25) now shares.
$values = @(2,2,2,1,1,1,2,2,2,1,2,2,2,1,2,2)

for ($r = 2; $r -le 25; $r++) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2
        $ws.Cells.Item($r, $col).Value = $values[$i]
    }
}
